$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.213.36'
$ws.Range("E2").Value = '  -2.83%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.821.38'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.11'
$ws.Range("E5").Value = '  -3.85%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.88'
$ws.Range("E6").Value = '  -4.34%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.816.71'
$ws.Range("E7").Value = '  +1.79%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("E9").Value = '  -0.25%  '

# Row 10
$ws.Range("E10").Value = '  -4.07%  '

# Row 11
$ws.Range("E11").Value = '  -0.47%  '

# Row 12
$ws.Range("E12").Value = '  -2.33%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.46'
$ws.Range("E13").Value = '  -4.29%  '

# Row 14
$ws.Range("E14").Value = '  -3.97%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.452.11'
$ws.Range("E15").Value = '  +1.88%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.819.14'
$ws.Range("E16").Value = '  +2.01%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.295.69'
$ws.Range("E17").Value = '  -2.73%  '

# Row 18
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.27'
$ws.Range("E18").Value = '  -4.24%  '

# Row 19
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.116'
$ws.Range("E19").Value = '  -4.49%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.06'
$ws.Range("E20").Value = '  -2.16%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '491.14'
$ws.Range("E21").Value = '  -2.86%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.36'
$ws.Range("E22").Value = '  +1.53%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.739'
$ws.Range("E23").Value = '  +2.26%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.11'
$ws.Range("E24").Value = '  -1.05%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  -7.33%  '

# Row 26
$ws.Range("E26").Value = '  +4.89%  '

# Row 27
$ws.Range("E27").Value = '  -5.71%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.19'
$ws.Range("E28").Value = '  -8.56%  '

# Row 29
$ws.Range("E29").Value = '  +0.02%  '

# Row 30
$ws.Range("E30").Value = '  +0.47%  '

# Row 31
$ws.Range("E31").Value = '  -1.37%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.13'
$ws.Range("E32").Value = '  +8.39%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.61'
$ws.Range("E33").Value = '  -3.48%  '

# Row 34
$ws.Range("E34").Value = '  -3.65%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.03%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("E36").Value = '  -3.99%  '

# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.137'
$ws.Range("E37").Value = '  -1.29%  '

# Row 38
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.83'
$ws.Range("E38").Value = '  -4.55%  '

# Row 39
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.326'
$ws.Range("E39").Value = '  -6.41%  '

# Row 40
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '452.96'
$ws.Range("E40").Value = '  +5.47%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '49.19'
$ws.Range("E41").Value = '  -1.39%  '

# Row 42
$ws.Range("E42").Value = '  -3.33%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.91'
$ws.Range("E43").Value = '  -9.44%  '

# Row 44
$ws.Range("E44").Value = '  -3.20%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.41'
$ws.Range("E45").Value = '  -7.12%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.855.63'
$ws.Range("E46").Value = '  -3.78%  '

# Row 47
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0353'
$ws.Range("E47").Value = '  -2.44%  '

# Row 48
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.04%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '137.93'
$ws.Range("E49").Value = '  +1.16%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.49'
$ws.Range("E50").Value = '  -2.88%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.70'
$ws.Range("E51").Value = '  +8.54%  '
